$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.423.05"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "1.849.31"
$ws.Range("D4").Value = "'0.9986"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'240.90"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").Value = "'0.6336"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").Value = "'0.9996"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.500.56"
$ws.Range("E8").Value = "  +89.12%  "
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").Value = "'0.2973"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").Value = "'24.66"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").Value = "3.755.62"
$ws.Range("E12").Value = "  +79.41%  "
$ws.Range("D13").Value = "'0.07711"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").Value = "'4.992"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").Value = "'0.6857"
$ws.Range("D16").Value = "'83.13"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "'0.000009968"
$ws.Range("E17").Value = "  +4.44%  "
$ws.Range("D18").Value = "'6.181"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").Value = "29.437.48"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "'232.07"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "'12.51"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "'7.614"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "'155.13"
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").Value = "'0.1386"
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("D27").Value = "'8.419"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "'1.469"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("D30").Value = "'0.05815"
$ws.Range("E30").Value = "  -2.95%  "
$ws.Range("B31").Value = "RocketPoolETH"
$ws.Range("C31").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D31").Value = "3.724.58"
$ws.Range("E31").Value = "  +85.41%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'1.259"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").Value = "'4.135"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("D35").Value = "'1.860"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("D37").Value = "'0.7170"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").Value = "1.251.53"
$ws.Range("E39").Value = "  +4.65%  "
$ws.Range("D40").Value = "'2.798"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").Value = "'0.01805"
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("D42").Value = "'0.8988"
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("D43").Value = "'6.098"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("D44").Value = "'0.9994"
$ws.Range("D45").Value = "'101.76"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("D47").Value = "'7.199"
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("D48").Value = "'9.170"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").Value = "'0.4019"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").Value = "'1.688"
$ws.Range("E50").Value = "  +1.86%  "
$ws.Range("D51").Value = "'0.1126"
$ws.Range("E51").Value = "  +0.25%  "
